# Weekly update: insert a new price-report row above row 172 (most recent
# week's data), pushing the existing rows 172:245 down to 173:246.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(172).Insert()

$ws.Cells.Item(172, 1).Value  = 7
$ws.Cells.Item(172, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value  = "Ñuble"
$ws.Cells.Item(172, 4).Value  = 44839
$ws.Cells.Item(172, 5).Value  = 16
$ws.Cells.Item(172, 6).Value  = 100112017
$ws.Cells.Item(172, 7).Value  = "Apio"
$ws.Cells.Item(172, 8).Value  = "Americana (o)"
$ws.Cells.Item(172, 9).Value  = "Primera"
$ws.Cells.Item(172, 10).Value = 120
$ws.Cells.Item(172, 11).Value = 9000
$ws.Cells.Item(172, 12).Value = 10000
$ws.Cells.Item(172, 13).Value = 9500
$ws.Cells.Item(172, 14).Value = "`$/docena de matas"
$ws.Cells.Item(172, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(172, 16).Value = 1583
$ws.Cells.Item(172, 17).Value = 6
$ws.Cells.Item(172, 18).Value = "Hortaliza"
